# DoFinance: switch to English website
#
# The DoFinance scraper now also reports deposit/withdrawal totals, so two
# new columns - "Einzahlungen" (deposits) and "Auszahlungen" (withdrawals) -
# are inserted right after the existing "Endsaldo" (closing balance) column
# on every sheet, between "Endsaldo" and "Investitionen". All following
# columns shift two places to the right. The new columns are filled with 0
# for every existing data row (matching the other numeric columns).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Tagesergebnisse" (daily results)
#   Header row 1: ... E=Endsaldo, F=Investitionen, ... -> insert at F:G
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Tagesergebnisse")

$ws1.Columns("F:G").Insert()
$ws1.Range("F1").Value = "Einzahlungen"
$ws1.Range("G1").Value = "Auszahlungen"
$ws1.Range("F1:G1").ColumnWidth = 15.15
$ws1.Range("F2").Value = 0
$ws1.Range("G2").Value = 0

# The trailing, entirely empty rows at the bottom of the original sheet are
# no longer written out.
$ws1.Rows("1048574:1048576").Delete()

$ws1.Rows(2).RowHeight = 15
$ws1.Range("A1").Select()

# ---------------------------------------------------------------------
# Sheet 2: "Monatsergebnisse" (monthly results)
#   Header row 1: ... E=Endsaldo, F=Investitionen, ... -> insert at F:G
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Monatsergebnisse")

$ws2.Columns("F:G").Insert()
$ws2.Range("F1").Value = "Einzahlungen"
$ws2.Range("G1").Value = "Auszahlungen"
$ws2.Range("F1:G1").ColumnWidth = 15.15
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 0
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 0
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = 0
$ws2.Range("F5").Value = 0
$ws2.Range("G5").Value = 0

$ws2.Rows(2).RowHeight = 15
$ws2.Range("A1").Select()

# ---------------------------------------------------------------------
# Sheet 3: "Gesamtergebnis" (overall result)
#   Header row 1: ... D=Endsaldo, E=Investitionen, ... -> insert at E:F
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Gesamtergebnis")

$ws3.Columns("E:F").Insert()
$ws3.Range("E1").Value = "Einzahlungen"
$ws3.Range("F1").Value = "Auszahlungen"
$ws3.Range("E1:F1").ColumnWidth = 15.15
$ws3.Range("E2").Value = 0
$ws3.Range("F2").Value = 0
$ws3.Range("E3").Value = 0
$ws3.Range("F3").Value = 0

$ws3.Rows(3).RowHeight = 15
$ws3.Range("A1").Select()

$ws1.Activate()
